# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45177 to 45178 (one day later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2 through 203 (row 1 is the header row).
$firstRow = 2
$lastRow = 203

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45178
